$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3079.25
$ws.Range("J2").Value = 7661.6665
$ws.Range("L2").Value = 7661.6665
$ws.Range("N2").Value = -7887.6665

$ws.Range("H12").Value = 401
$ws.Range("I12").Value = 401
$ws.Range("K12").Value = 401
$ws.Range("M12").Value = -231

$ws.Range("H86").Value = 1659.8
$ws.Range("I86").Value = 1477.5555
$ws.Range("K86").Value = 1477.5555
$ws.Range("M86").Value = -354.5554999999999

$ws.Range("H89").Value = 1659.8
$ws.Range("I89").Value = 1477.5555
$ws.Range("K89").Value = 7387.7775
$ws.Range("M89").Value = -1771.7775

$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 400
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 848
$ws.Range("N92").ClearContents()

$ws.Range("H99").Value = 200.33333
$ws.Range("I99").Value = 200.33333
$ws.Range("K99").Value = 600.99999
$ws.Range("M99").Value = 897.00001

$ws.Range("H101").Value = 16670146
$ws.Range("J101").Value = 993
$ws.Range("L101").Value = 2979
$ws.Range("N101").Value = -6223

$ws.Range("H111").Value = 10490.823
$ws.Range("I111").Value = 12671.091
$ws.Range("K111").Value = 38013.273
$ws.Range("M111").Value = -34946.273

$ws.Range("H116").Value = 5931.6665
$ws.Range("I116").Value = 4497.5
$ws.Range("J116").Value = 8800
$ws.Range("K116").Value = 4497.5
$ws.Range("L116").Value = 8800
$ws.Range("M116").Value = -1055.5
$ws.Range("N116").Value = -15684

$ws.Range("H118").Value = 1120
$ws.Range("I118").Value = 1120
$ws.Range("K118").Value = 3360
$ws.Range("M118").Value = -1703

$ws.Range("H125").Value = 8614.333000000001
$ws.Range("I125").Value = 1037.25
$ws.Range("K125").Value = 9335.25
$ws.Range("M125").Value = -6875.25

$ws.Range("H132").Value = 3942.8
$ws.Range("I132").Value = 3942.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11828.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9298.400000000001
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 5799.1
$ws.Range("I137").Value = 4166.3335
$ws.Range("K137").Value = 12499.0005
$ws.Range("M137").Value = -9949.000499999998

$ws.Range("H138").Value = 2501.15
$ws.Range("I138").Value = 1840
$ws.Range("J138").Value = 3309.2222
$ws.Range("K138").Value = 5520
$ws.Range("L138").Value = 9927.6666
$ws.Range("M138").Value = -380
$ws.Range("N138").Value = -20207.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7425
$ws.Range("I45").Value = 13000
$ws.Range("J45").Value = 1850
$ws.Range("K45").Value = 13000
$ws.Range("L45").Value = 1850
$ws.Range("M45").Value = -12623
$ws.Range("N45").Value = -2604

$ws.Range("H119").Value = 32832.332
$ws.Range("J119").Value = 32832.332
$ws.Range("L119").Value = 32832.332
$ws.Range("N119").Value = -42508.332

$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820

$ws.Range("H132").Value = 1749
$ws.Range("I132").Value = 1749
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5247
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2717
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2147.1
$ws.Range("I99").Value = 1492.8572
$ws.Range("J99").Value = 3673.6667
$ws.Range("K99").Value = 1492.8572
$ws.Range("L99").Value = 3673.6667
$ws.Range("M99").Value = 5.142800000000079
$ws.Range("N99").Value = -6669.6667

$ws.Range("H132").Value = 75999
$ws.Range("J132").Value = 75999
$ws.Range("L132").Value = 75999
$ws.Range("N132").Value = -86119

$ws.Range("H134").Value = 1416.6666
$ws.Range("J134").Value = 1000
$ws.Range("L134").Value = 3000
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2670.7273
$ws.Range("I31").Value = 2098.3333
$ws.Range("K31").Value = 2098.3333
$ws.Range("M31").Value = -1803.3333

$ws.Range("H34").Value = 2670.7273
$ws.Range("I34").Value = 2098.3333
$ws.Range("K34").Value = 2098.3333
$ws.Range("M34").Value = -1896.3333

$ws.Range("H99").Value = 1645.6428
$ws.Range("I99").Value = 1398.909
$ws.Range("J99").Value = 2550.3333
$ws.Range("K99").Value = 1398.909
$ws.Range("L99").Value = 2550.3333
$ws.Range("M99").Value = 99.09099999999989
$ws.Range("N99").Value = -5546.3333

$ws.Range("H122").Value = 1624.6666
$ws.Range("I122").Value = 1624.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4873.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2423.9998
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 1645.6428
$ws.Range("I126").Value = 1398.909
$ws.Range("J126").Value = 2550.3333
$ws.Range("K126").Value = 4196.727000000001
$ws.Range("L126").Value = 7650.999899999999
$ws.Range("M126").Value = -1726.727000000001
$ws.Range("N126").Value = -12590.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 183.5
$ws.Range("I14").Value = 183.5
$ws.Range("K14").Value = 550.5
$ws.Range("M14").Value = -377.5

$ws.Range("H128").Value = 542941
$ws.Range("I128").Value = 542941
$ws.Range("K128").Value = 1628823
$ws.Range("M128").Value = -1623843

$ws.Range("H131").Value = 668359.4
$ws.Range("J131").Value = 1431003.6
$ws.Range("L131").Value = 4293010.800000001
$ws.Range("N131").Value = -4303090.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1444.4286
$ws.Range("J80").Value = 1102.75
$ws.Range("L80").Value = 1102.75
$ws.Range("N80").Value = -3098.75

$ws.Range("H83").Value = 1444.4286
$ws.Range("J83").Value = 1102.75
$ws.Range("L83").Value = 5513.75
$ws.Range("N83").Value = -15497.75

$ws.Range("H97").Value = 197.77777
$ws.Range("I97").Value = 147.14285
$ws.Range("K97").Value = 147.14285
$ws.Range("M97").Value = 348.85715

$ws.Range("H134").Value = 47665.2
$ws.Range("J134").Value = 47665.2
$ws.Range("L134").Value = 142995.6
$ws.Range("N134").Value = -148065.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9321.647000000001
$ws.Range("I7").Value = 8746.333000000001
$ws.Range("K7").Value = 8746.333000000001
$ws.Range("M7").Value = -8634.333000000001

$ws.Range("H16").Value = 1634.2858
$ws.Range("I16").Value = 1790.75
$ws.Range("J16").Value = 1425.6666
$ws.Range("K16").Value = 1790.75
$ws.Range("L16").Value = 1425.6666
$ws.Range("M16").Value = -1620.75
$ws.Range("N16").Value = -1765.6666

$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 3000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2705

$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 3000
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = -2893

$ws.Range("H55").Value = 1051
$ws.Range("I55").Value = 603.6667
$ws.Range("J55").Value = 1386.5
$ws.Range("K55").Value = 603.6667
$ws.Range("L55").Value = 1386.5
$ws.Range("M55").Value = -430.6667
$ws.Range("N55").Value = -1732.5

$ws.Range("H126").Value = 9321.647000000001
$ws.Range("I126").Value = 8746.333000000001
$ws.Range("K126").Value = 26238.999
$ws.Range("M126").Value = -23768.999

$ws.Range("H132").Value = 2849.818
$ws.Range("I132").Value = 2574
$ws.Range("J132").Value = 3180.8
$ws.Range("K132").Value = 7722
$ws.Range("L132").Value = 9542.400000000001
$ws.Range("M132").Value = -5192
$ws.Range("N132").Value = -14602.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws.Range("H136").Value = 1654.9
$ws.Range("I136").Value = 1142.7142
$ws.Range("K136").Value = 3428.1426
$ws.Range("M136").Value = -878.1425999999997
